$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-13 20:18:41"
$ws.Range("I2").Value = "2.3 mm"
$ws.Range("E3").Value = "2026-02-13 20:18:44"
$ws.Range("I3").Value = "6.1 mm"
$ws.Range("E4").Value = "2026-02-13 20:18:47"
$ws.Range("J4").Value = "994.2 hPa"
$ws.Range("E5").Value = "2026-02-13 20:18:49"
$ws.Range("G5").Value = "108 cm"
$ws.Range("I5").Value = "1.3 mm"
$ws.Range("E6").Value = "2026-02-13 20:18:52"
$ws.Range("H6").Value = "'77%"
$ws.Range("J6").Value = "994.3 hPa"
$ws.Range("E7").Value = "2026-02-13 20:18:55"
$ws.Range("J7").Value = "994.6 hPa"
$ws.Range("O7").Value = "12.8 °C"
$ws.Range("E8").Value = "2026-02-13 20:18:57"
$ws.Range("J8").Value = "994.5 hPa"
$ws.Range("N8").Value = "6.8 °C 19:57 TU"
$ws.Range("O8").Value = "9.2 °C"
$ws.Range("E9").Value = "2026-02-13 20:19:00"
$ws.Range("L9").Value = "33.8 km/h - 332º 19:51 TU"
$ws.Range("E10").Value = "2026-02-13 20:19:03"
$ws.Range("I10").Value = "20.2 mm"
$ws.Range("E11").Value = "2026-02-13 20:19:06"
$ws.Range("E12").Value = "2026-02-13 20:19:08"
$ws.Range("O12").Value = "9.4 °C"
$ws.Range("E13").Value = "2026-02-13 20:19:11"
$ws.Range("E14").Value = "2026-02-13 20:19:14"
$ws.Range("H14").Value = "'85%"
$ws.Range("O14").Value = "10.4 °C"
$ws.Range("E15").Value = "2026-02-13 20:19:16"
$ws.Range("I15").Value = "4.8 mm"
$ws.Range("E16").Value = "2026-02-13 20:19:19"
$ws.Range("I16").Value = "12.9 mm"
$ws.Range("E17").Value = "2026-02-13 20:19:21"
$ws.Range("G17").Value = "1 cm"
$ws.Range("E18").Value = "2026-02-13 20:19:24"
$ws.Range("J18").Value = "994.4 hPa"
$ws.Range("E19").Value = "2026-02-13 20:19:27"
$ws.Range("E20").Value = "2026-02-13 20:19:29"
$ws.Range("I20").Value = "23.1 mm"
$ws.Range("E21").Value = "2026-02-13 20:19:32"
$ws.Range("H21").Value = "'92%"
$ws.Range("J21").Value = "997.3 hPa"
$ws.Range("E22").Value = "2026-02-13 20:19:35"
$ws.Range("E23").Value = "2026-02-13 20:19:37"
$ws.Range("I23").Value = "10.1 mm"
$ws.Range("E24").Value = "2026-02-13 20:19:40"
$ws.Range("H24").Value = "'96%"
$ws.Range("J24").Value = "995.3 hPa"
$ws.Range("E25").Value = "2026-02-13 20:19:43"
$ws.Range("G25").Value = "116 cm"
$ws.Range("I25").Value = "9.1 mm"
$ws.Range("E26").Value = "2026-02-13 20:19:45"
$ws.Range("E27").Value = "2026-02-13 20:19:48"
$ws.Range("E28").Value = "2026-02-13 20:19:50"
$ws.Range("H28").Value = "'80%"
$ws.Range("J28").Value = "994.7 hPa"
$ws.Range("E29").Value = "2026-02-13 20:19:53"
$ws.Range("I29").Value = "14.0 mm"
$ws.Range("L29").Value = "29.9 km/h - 329º 19:43 TU"
$ws.Range("E30").Value = "2026-02-13 20:19:55"
$ws.Range("J30").Value = "994.2 hPa"
$ws.Range("E31").Value = "2026-02-13 20:19:58"
$ws.Range("I31").Value = "4.7 mm"
$ws.Range("J31").Value = "993.1 hPa"
$ws.Range("E32").Value = "2026-02-13 20:20:01"
$ws.Range("E33").Value = "2026-02-13 20:20:03"
$ws.Range("H33").Value = "'90%"
$ws.Range("J33").Value = "996.2 hPa"
$ws.Range("E34").Value = "2026-02-13 20:20:06"
$ws.Range("E35").Value = "2026-02-13 20:20:08"
$ws.Range("I35").Value = "8.1 mm"
$ws.Range("J35").Value = "995.3 hPa"
$ws.Range("E36").Value = "2026-02-13 20:20:11"
$ws.Range("I36").Value = "8.9 mm"
$ws.Range("J36").Value = "994.3 hPa"
$ws.Range("L36").Value = "49.7 km/h - 325º 19:49 TU"
$ws.Range("E37").Value = "2026-02-13 20:20:13"
$ws.Range("H37").Value = "'85%"
$ws.Range("J37").Value = "996.2 hPa"
$ws.Range("E38").Value = "2026-02-13 20:20:16"
$ws.Range("E39").Value = "2026-02-13 20:20:18"
$ws.Range("I39").Value = "19.4 mm"
$ws.Range("E40").Value = "2026-02-13 20:20:21"
$ws.Range("H40").Value = "'98%"
$ws.Range("J40").Value = "997.8 hPa"
$ws.Range("E41").Value = "2026-02-13 20:20:23"
$ws.Range("H41").Value = "'76%"
$ws.Range("J41").Value = "994.7 hPa"
$ws.Range("K41").Value = "1.5 MJ/m2"
$ws.Range("E42").Value = "2026-02-13 20:20:26"
$ws.Range("I42").Value = "10.5 mm"
$ws.Range("E43").Value = "2026-02-13 20:20:28"
$ws.Range("E44").Value = "2026-02-13 20:20:31"
$ws.Range("I44").Value = "7.5 mm"
$ws.Range("O44").Value = "-4.0 °C"
$ws.Range("E45").Value = "2026-02-13 20:20:33"
$ws.Range("H45").Value = "'65%"
$ws.Range("I45").Value = "0.5 mm"
$ws.Range("J45").Value = "993.2 hPa"
$ws.Range("E46").Value = "2026-02-13 20:20:36"
$ws.Range("J46").Value = "995.4 hPa"
